# Auto-generated edit script: updates crypto price/volume data to the
# latest GitHub Actions scrape, matching the authored commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.943.55"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.644.38"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'213.43"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").Value = "'0.5219"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.2606"
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").Value = "'0.06334"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "'20.67"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").Value = "'0.07682"
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").Value = "1.642.51"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "'4.418"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "1.864.56"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "'0.5495"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "0.0₅8213"
$ws.Range("E16").Value = "  +3.57%  "
$ws.Range("D17").Value = "'64.59"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "25.925.37"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'4.700"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'189.35"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "'10.18"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'6.265"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "'142.88"
$ws.Range("E25").Value = "  -4.00%  "
$ws.Range("D26").Value = "'0.1242"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").Value = "'7.378"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "'15.93"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").Value = "'1.404"
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("D30").Value = "'0.05938"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").Value = "'1.256"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'3.413"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.402"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").Value = "'1.639"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "'0.9902"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'2.393"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'2.741"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").Value = "'0.5614"
$ws.Range("E38").Value = "  -4.49%  "
$ws.Range("D39").Value = "'0.01603"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").Value = "'5.844"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").Value = "'0.8546"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "1.024.29"
$ws.Range("E43").Value = "  -7.30%  "
$ws.Range("D44").Value = "'98.76"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").Value = "1.788.92"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'55.57"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.021"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05137"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.4208"
$ws.Range("E51").Value = "  -0.94%  "
